$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the two new header cells, copying the header style (bold, centered, bordered) from A1
$ws.Range("A1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$ws.Range("E1").Value = "reviews_count"
$ws.Range("F1").Value = "reviews_average"

# Row 2
$ws.Range("A2").Value = "Dr Pascal MARIN"
$ws.Range("B2").Value = "5 Rue Crétet, 75009 Paris, France"
$ws.Range("C2").Value = "doctolib.fr"
$ws.Range("D2").Value = "+33 6 75 15 49 16"
$ws.Range("E2").Value = 5
$ws.Range("F2").Value = 3.4

# Row 3
$ws.Range("A3").Value = "Dr. Charlotte Parment"
$ws.Range("B3").Value = "cabinet médical ipso Saint Martin, 323 Rue Saint-Martin, 75003 Paris, France"
$ws.Range("C3").Value = "ipso.paris"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 5

# Row 4
$ws.Range("A4").Value = "Dr Claire Paris"
$ws.Range("B4").Value = "86 Rue de l'Université, 75007 Paris, France"
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = "+33 1 40 62 95 28"
$ws.Range("E4").Value = 6
$ws.Range("F4").Value = 5

# Row 5
$ws.Range("A5").Value = "Dr Nancy Salzman"
$ws.Range("B5").Value = "1 Av. de Lowendal, 75007 Paris, France"
$ws.Range("C5").Value = "doctor-salzman.com"
$ws.Range("D5").Value = "+33 1 45 63 18 43"
$ws.Range("E5").Value = 35
$ws.Range("F5").Value = 4.8

# Row 6
$ws.Range("A6").Value = "Docteur Franck Besse"
$ws.Range("B6").Value = "45 Rue de Lancry, 75010 Paris, France"
$ws.Range("C6").Value = "doctolib.fr"
$ws.Range("D6").Value = "+33 1 44 85 26 83"
$ws.Range("E6").Value = 24
$ws.Range("F6").Value = 4.2

# Row 7 (new row)
$ws.Range("A7").Value = "Docteur Simon OHAYON- English speaking doctor- International medical center"
$ws.Range("B7").Value = "48 BIS Rue des Belles Feuilles, 75116 Paris, France"
$ws.Range("C7").Value = "doctolib.fr"
$ws.Range("D7").Value = "+33 6 58 80 18 38"
$ws.Range("E7").Value = 94
$ws.Range("F7").Value = 4.1
